$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 ("I0") and J1 ("IF") matching the style used by the
# other header cells (bold font, thin border all around, centered
# horizontal/top vertical alignment).
$ws.Range("I1").Value = "I0"
$ws.Range("I1").Font.Bold = $true
$ws.Range("I1").HorizontalAlignment = -4108
$ws.Range("I1").VerticalAlignment = -4160
$ws.Range("I1").Borders.LineStyle = 1

$ws.Range("J1").Value = "IF"
$ws.Range("J1").Font.Bold = $true
$ws.Range("J1").HorizontalAlignment = -4108
$ws.Range("J1").VerticalAlignment = -4160
$ws.Range("J1").Borders.LineStyle = 1

# Data rows 2..31 for columns I (I0) and J (IF)
$iValues = @(1,1,1,6,6,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,5,6)
$jValues = @(3,5,3,8,6,5,5,5,7,5,5,4,5,4,4,5,5,5,6,4,5,6,5,3,6,4,5,4,7,7)

for ($r = 0; $r -lt 30; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$r]
    $ws.Cells.Item($row, 10).Value = $jValues[$r]
}
